# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed quote data, matching a GitHub Actions scraper
# run. Numeric-looking Price values are forced to remain plain text
# (matching the sheet's original string-typed cells) by temporarily
# switching the cell to a Text number format before assignment, then
# resetting the style back to Normal so no stray style index is left
# behind in the saved file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.831.57'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '2.956.90'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '2.955.85'
$ws.Range('E8').Value = '  +1.01%  '
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.444'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000239'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.16%  '
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '3.451.13'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '62.725.81'
$ws.Range('E17').Value = '  +2.80%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '2.959.84'
$ws.Range('E19').Value = '  +1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '441.48'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  -1.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('E26').Value = '  -3.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.45%  '
$ws.Range('E30').Value = '  +0.44%  '
$ws.Range('E31').Value = '  -3.42%  '
$ws.Range('D32').Value = '0.0₃0975'
$ws.Range('E32').Value = '  +10.84%  '
$ws.Range('E33').Value = '  -0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.05'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.52'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('E42').Value = '  -4.50%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').Value = '2.736.17'
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '366.15'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0339'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.93'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.68%  '
